$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C6').Value = 'Remis'
$ws.Range('C8').Value = 'Radomiak Radom'
$ws.Range('C9').Value = 'Warta Poznań'
$ws.Range('C10').Value = 'Lechia Gdańsk'
$ws.Range('C12').Value = 'Raków Częstochowa'
$ws.Range('C14').Value = 'Remis'
$ws.Range('C17').Value = 'Remis'
$ws.Range('C20').Value = 'Legia Warszawa'
$ws.Range('C21').Value = 'Zagłębie Lubin'
$ws.Range('C24').Value = 'Warta Poznań'
$ws.Range('C26').Value = 'Remis'
$ws.Range('C27').Value = 'Raków Częstochowa'
$ws.Range('C28').Value = 'Lechia Gdańsk'
$ws.Range('C29').Value = 'Radomiak Radom'
$ws.Range('C30').Value = 'Lech Poznań'
$ws.Range('C32').Value = 'Remis'
$ws.Range('C37').Value = 'Remis'
$ws.Range('C42').Value = 'Remis'
$ws.Range('C44').Value = 'Lechia Gdańsk'
$ws.Range('C45').Value = 'Remis'
$ws.Range('C49').Value = 'Radomiak Radom'
$ws.Range('C50').Value = 'Remis'
$ws.Range('C52').Value = 'Pogoń Szczecin'
$ws.Range('C53').Value = 'Widzew Łódź'
$ws.Range('C54').Value = 'Remis'
$ws.Range('C55').Value = 'Remis'
$ws.Range('C56').Value = 'Cracovia'
$ws.Range('C58').Value = 'Lechia Gdańsk'
$ws.Range('C61').Value = 'Radomiak Radom'
$ws.Range('C63').Value = 'Remis'
$ws.Range('C65').Value = 'Raków Częstochowa'
$ws.Range('C66').Value = 'Miedź Legnica'
$ws.Range('C67').Value = 'Zagłębie Lubin'
$ws.Range('C72').Value = 'Śląsk Wrocław'
$ws.Range('C76').Value = 'Piast Gliwice'
$ws.Range('C78').Value = 'Remis'
$ws.Range('C79').Value = 'Legia Warszawa'
$ws.Range('C80').Value = 'Remis'
$ws.Range('C81').Value = 'Cracovia'
$ws.Range('C82').Value = 'Lechia Gdańsk'
$ws.Range('C84').Value = 'Piast Gliwice'
$ws.Range('C85').Value = 'Zagłębie Lubin'
$ws.Range('C86').Value = 'Korona Kielce'
$ws.Range('C94').Value = 'Miedź Legnica'
$ws.Range('C96').Value = 'Lechia Gdańsk'
$ws.Range('C97').Value = 'Radomiak Radom'
$ws.Range('C99').Value = 'Wisła Płock'
$ws.Range('C100').Value = 'Warta Poznań'
$ws.Range('C101').Value = 'Remis'
$ws.Range('C102').Value = 'Widzew Łódź'
$ws.Range('C107').Value = 'Miedź Legnica'
$ws.Range('C108').Value = 'Remis'
$ws.Range('C109').Value = 'Remis'
$ws.Range('C110').Value = 'Korona Kielce'
$ws.Range('C112').Value = 'Cracovia'
$ws.Range('C114').Value = 'Pogoń Szczecin'
$ws.Range('C115').Value = 'Radomiak Radom'
$ws.Range('C117').Value = 'Zagłębie Lubin'
$ws.Range('C118').Value = 'Wisła Płock'
$ws.Range('C121').Value = 'Zagłębie Lubin'
$ws.Range('C123').Value = 'Remis'
$ws.Range('C124').Value = 'Stal Mielec'
$ws.Range('C125').Value = 'Górnik Zabrze'
$ws.Range('C126').Value = 'Miedź Legnica'
$ws.Range('C128').Value = 'Jagielonia Białystok'
$ws.Range('C129').Value = 'Cracovia'
$ws.Range('C130').Value = 'Remis'
$ws.Range('C131').Value = 'Lech Poznań'
$ws.Range('C132').Value = 'Remis'
$ws.Range('C137').Value = 'Cracovia'
$ws.Range('C138').Value = 'Piast Gliwice'
$ws.Range('C140').Value = 'Legia Warszawa'
$ws.Range('C141').Value = 'Śląsk Wrocław'
$ws.Range('C142').Value = 'Remis'
$ws.Range('C145').Value = 'Remis'
$ws.Range('C146').Value = 'Remis'
$ws.Range('C147').Value = 'Raków Częstochowa'
$ws.Range('C149').Value = 'Lechia Gdańsk'
$ws.Range('C150').Value = 'Górnik Zabrze'
$ws.Range('C156').Value = 'Jagielonia Białystok'
$ws.Range('C158').Value = 'Korona Kielce'
$ws.Range('C163').Value = 'Remis'
$ws.Range('C164').Value = 'Legia Warszawa'
$ws.Range('C165').Value = 'Cracovia'
$ws.Range('C167').Value = 'Remis'
$ws.Range('C170').Value = 'Raków Częstochowa'
$ws.Range('C171').Value = 'Remis'
$ws.Range('C172').Value = 'Remis'
$ws.Range('C176').Value = 'Legia Warszawa'
$ws.Range('C177').Value = 'Górnik Zabrze'
$ws.Range('C182').Value = 'Remis'
$ws.Range('C183').Value = 'Legia Warszawa'
$ws.Range('C185').Value = 'Remis'
$ws.Range('C188').Value = 'Jagielonia Białystok'
$ws.Range('C193').Value = 'Miedź Legnica'
$ws.Range('C194').Value = 'Lechia Gdańsk'
$ws.Range('C197').Value = 'Warta Poznań'
$ws.Range('C198').Value = 'Wisła Płock'
$ws.Range('C200').Value = 'Śląsk Wrocław'
$ws.Range('C201').Value = 'Remis'
$ws.Range('C203').Value = 'Remis'
$ws.Range('C204').Value = 'Górnik Zabrze'
$ws.Range('C206').Value = 'Remis'
$ws.Range('C207').Value = 'Remis'
$ws.Range('C209').Value = 'Lech Poznań'
$ws.Range('C212').Value = 'Radomiak Radom'
$ws.Range('C216').Value = 'Remis'
$ws.Range('C218').Value = 'Jagielonia Białystok'
$ws.Range('C219').Value = 'Piast Gliwice'
$ws.Range('C220').Value = 'Wisła Płock'
$ws.Range('C227').Value = 'Widzew Łódź'
$ws.Range('C229').Value = 'Warta Poznań'
$ws.Range('C230').Value = 'Remis'
$ws.Range('C231').Value = 'Remis'
$ws.Range('C235').Value = 'Radomiak Radom'
$ws.Range('C236').Value = 'Lechia Gdańsk'
$ws.Range('C242').Value = 'Remis'
$ws.Range('C243').Value = 'Remis'
$ws.Range('C244').Value = 'Śląsk Wrocław'
$ws.Range('C245').Value = 'Remis'
$ws.Range('C249').Value = 'Remis'
$ws.Range('C250').Value = 'Remis'
$ws.Range('C252').Value = 'Remis'
$ws.Range('C253').Value = 'Warta Poznań'
$ws.Range('C254').Value = 'Remis'
$ws.Range('C255').Value = 'Lechia Gdańsk'
$ws.Range('C256').Value = 'Raków Częstochowa'
$ws.Range('C257').Value = 'Remis'
$ws.Range('C258').Value = 'Remis'
$ws.Range('C260').Value = 'Remis'
$ws.Range('C261').Value = 'Widzew Łódź'
$ws.Range('C264').Value = 'Pogoń Szczecin'
$ws.Range('C269').Value = 'Lechia Gdańsk'
$ws.Range('C270').Value = 'Remis'
$ws.Range('C274').Value = 'Remis'
$ws.Range('C276').Value = 'Remis'
$ws.Range('C277').Value = 'Remis'
$ws.Range('C278').Value = 'Pogoń Szczecin'
$ws.Range('C280').Value = 'Remis'
$ws.Range('C287').Value = 'Warta Poznań'
$ws.Range('C288').Value = 'Górnik Zabrze'
$ws.Range('C289').Value = 'Śląsk Wrocław'
$ws.Range('C291').Value = 'Stal Mielec'
$ws.Range('C293').Value = 'Legia Warszawa'
$ws.Range('C294').Value = 'Pogoń Szczecin'
$ws.Range('C298').Value = 'Remis'
$ws.Range('C299').Value = 'Remis'
$ws.Range('C302').Value = 'Legia Warszawa'
$ws.Range('C303').Value = 'Górnik Zabrze'
$ws.Range('C305').Value = 'Remis'
$ws.Range('C307').Value = 'Remis'

Write-Output "Updated 151 match result cells"